$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.023.00'
$ws.Range("E2").Value = '  +6.90%  '

$ws.Range("D3").Value = '3.565.57'
$ws.Range("E3").Value = '  +2.75%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '416.44'
$ws.Range("E5").Value = '  +0.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.58'
$ws.Range("E6").Value = '  -0.37%  '

$ws.Range("E7").Value = '  +3.45%  '

$ws.Range("D8").Value = '3.558.52'
$ws.Range("E8").Value = '  +2.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.768'
$ws.Range("E10").Value = '  +5.70%  '

$ws.Range("E11").Value = '  +12.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000330'
$ws.Range("E12").Value = '  +45.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.25'
$ws.Range("E13").Value = '  -0.72%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.93'
$ws.Range("E14").Value = '  +1.50%  '

$ws.Range("D15").Value = '4.139.26'
$ws.Range("E15").Value = '  +2.85%  '

$ws.Range("E16").Value = '  -0.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.40'
$ws.Range("E17").Value = '  -0.92%  '

$ws.Range("D18").Value = '3.593.61'
$ws.Range("E18").Value = '  +3.97%  '

$ws.Range("E19").Value = '  +5.08%  '

$ws.Range("D20").Value = '66.922.62'
$ws.Range("E20").Value = '  +6.77%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.33'
$ws.Range("E21").Value = '  -2.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '456.02'
$ws.Range("E22").Value = '  -1.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '88.13'
$ws.Range("E23").Value = '  -2.61%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.11'
$ws.Range("E24").Value = '  -5.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.42'
$ws.Range("E25").Value = '  +1.20%  '

$ws.Range("E26").Value = '  +1.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.11'
$ws.Range("E27").Value = '  -6.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '34.64'
$ws.Range("E28").Value = '  +3.98%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.87'
$ws.Range("E29").Value = '  +1.54%  '

$ws.Range("E30").Value = '  +4.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.34'
$ws.Range("E31").Value = '  +1.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.118'
$ws.Range("E32").Value = '  +4.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.37'
$ws.Range("E33").Value = '  -3.11%  '

$ws.Range("E34").Value = '  -4.68%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '40.94'
$ws.Range("E35").Value = '  -0.19%  '

$ws.Range("E36").Value = '  -0.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.57'
$ws.Range("E37").Value = '  -2.68%  '

$ws.Range("E38").Value = '  +0.65%  '

$ws.Range("D39").Value = '0.0₃0730'
$ws.Range("E39").Value = '  +27.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.147'
$ws.Range("E40").Value = '  +8.96%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").Value = '  -0.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.06'
$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '149.22'
$ws.Range("E43").Value = '  -0.25%  '

$ws.Range("E44").Value = '  +0.33%  '

$ws.Range("E45").Value = '  -2.15%  '

$ws.Range("E46").Value = '  -3.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.32'
$ws.Range("E47").Value = '  -2.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.98'
$ws.Range("E48").Value = '  -4.09%  '

$ws.Range("E49").Value = '  -2.12%  '

$ws.Range("E50").Value = '  +14.66%  '

$ws.Range("E51").Value = '  -4.73%  '
